$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.157.94"
$ws.Range("E2").Value = "  -2.01%  "

$ws.Range("D3").Value = "2.976.09"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "594.25"
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").Value = "141.96"
$ws.Range("E6").Value = "  -3.14%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  -1.61%  "

$ws.Range("D9").Value = "2.972.53"
$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").Value = "5.94"
$ws.Range("E11").Value = "  +3.86%  "

$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +2.03%  "

$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").Value = "34.04"
$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("E15").Value = "  +1.85%  "

$ws.Range("D16").Value = "3.463.63"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").Value = "61.228.41"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").Value = "6.84"
$ws.Range("E18").Value = "  -2.62%  "

$ws.Range("D19").Value = "2.969.84"
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").Value = "448.57"
$ws.Range("E20").Value = "  -2.28%  "

$ws.Range("D21").Value = "'14.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("D22").Value = "0.678"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").Value = "7.25"
$ws.Range("E23").Value = "  -1.73%  "

$ws.Range("D24").Value = "82.29"
$ws.Range("E24").Value = "  +2.82%  "

$ws.Range("D25").Value = "2.15"
$ws.Range("E25").Value = "  -5.26%  "

$ws.Range("D26").Value = "11.92"
$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("D27").Value = "10.25"
$ws.Range("E27").Value = "  +2.21%  "

$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "2.66"
$ws.Range("E29").Value = "  +1.72%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("D32").Value = "2.03"
$ws.Range("E32").Value = "  -2.82%  "

$ws.Range("D33").Value = "27.12"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("E35").Value = "  +1.97%  "

$ws.Range("E36").Value = "  -1.56%  "

$ws.Range("D37").Value = "5.74"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "50.23"
$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("E39").Value = "  -3.35%  "

$ws.Range("D40").Value = "8.95"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  +7.46%  "

$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  -4.21%  "

$ws.Range("D43").Value = "'388.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.15%  "

$ws.Range("D44").Value = "0.0347"
$ws.Range("E44").Value = "  -1.71%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").Value = "  -4.25%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "38.36"
$ws.Range("E46").Value = "  -1.99%  "

$ws.Range("D47").Value = "2.686.50"
$ws.Range("E47").Value = "  -3.06%  "

$ws.Range("D48").Value = "129.51"
$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").Value = "2.12"
$ws.Range("E51").Value = "  -0.85%  "
